# Auto-generated script to update cryptos worksheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.137.42'
$ws.Range('E2').Value = '  +2.97%  '
$ws.Range('D3').Value = '2.312.91'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '310.20'
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.47'
$ws.Range('E6').Value = '  +6.02%  '
$ws.Range('E7').Value = '  +2.44%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.525'
$ws.Range('E9').Value = '  +7.40%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.16'
$ws.Range('E10').Value = '  +3.73%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0819'
$ws.Range('E11').Value = '  +3.76%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.27'
$ws.Range('E13').Value = '  +8.21%  '
$ws.Range('D14').Value = '2.668.75'
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.00'
$ws.Range('E15').Value = '  +4.65%  '
$ws.Range('D16').Value = '2.311.18'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.814'
$ws.Range('E17').Value = '  +3.25%  '
$ws.Range('D18').Value = '43.086.87'
$ws.Range('E18').Value = '  +3.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.58'
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  +2.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.11'
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '68.45'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '241.13'
$ws.Range('E23').Value = '  +1.91%  '
$ws.Range('E24').Value = '  +5.58%  '
$ws.Range('E25').Value = '  +2.75%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.67'
$ws.Range('E27').Value = '  +4.55%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.31'
$ws.Range('E28').Value = '  +2.65%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.66'
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '166.11'
$ws.Range('E31').Value = '  +3.63%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.34'
$ws.Range('E32').Value = '  +3.22%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').Value = '  -1.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '17.93'
$ws.Range('E35').Value = '  +5.22%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0743'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.108'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.84'
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.116'
$ws.Range('E40').Value = '  +2.25%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.31'
$ws.Range('E41').Value = '  +8.54%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.31'
$ws.Range('E42').Value = '  +1.13%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0290'
$ws.Range('E43').Value = '  +2.92%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.974.97'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '19.18'
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('E46').Value = '  +3.99%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.82'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.96'
$ws.Range('E48').Value = '  +18.44%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '55.64'
$ws.Range('E49').Value = '  +5.49%  '
$ws.Range('D50').Value = '2.536.68'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('E51').Value = '  +2.92%  '
